$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44319, 5, 30, 166.8706196462343),
    @(44320, 1, 24, 133.4964957169874),
    @(44321, 1, 24, 133.4964957169874)
)

$lastRow = 244
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the format of the last existing row down to the new row so the
    # date column keeps its date style/number format.
    $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 4)).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 4)).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = $false
